$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting from the last existing date cell so the new
# date cells reuse the existing style (s="1") instead of creating a new one.
$ws.Range("A11").Copy()
$ws.Range("A12:A13").PasteSpecial(-4122)

# Row 12: Volta Comunitat Valenciana - Stage 4
$ws.Range("A12").Value = "2/7/2026"
$ws.Range("B12").Value = "Volta Comunitat Valenciana"
$ws.Range("C12").Value = "Stage 4"
$ws.Range("D12").Value = "Remco Evenepoel"
$ws.Range("E12").Value = "João Almeida"
$ws.Range("F12").Value = "Giulio Pellizzari"
$ws.Range("G12").Value = "Antonio Tiberi"
$ws.Range("H12").Value = "Brandon McNulty"
$ws.Range("I12").Value = "Magnus Sheffield"
$ws.Range("J12").Value = "Aleksandr Vlasov"
$ws.Range("K12").Value = "Ben Turner"
$ws.Range("L12").Value = "Riley Sheehan"
$ws.Range("M12").Value = "Viktor Soenens"

# Row 13: Volta Comunitat Valenciana - Stage 5
$ws.Range("A13").Value = "2/8/2026"
$ws.Range("B13").Value = "Volta Comunitat Valenciana"
$ws.Range("C13").Value = "Stage 5"
$ws.Range("D13").Value = "Raúl García Pierna"
$ws.Range("E13").Value = "Emil Herzog"
$ws.Range("F13").Value = "Jasper Schoofs"
$ws.Range("G13").Value = "Adrià Pericas"
$ws.Range("H13").Value = "Sven Erik Bystrøm"
$ws.Range("I13").Value = "Diego Uriarte"
$ws.Range("J13").Value = "Ben Turner"
$ws.Range("K13").Value = "Dries Van Gestel"
$ws.Range("L13").Value = "Mathias Vacek"
$ws.Range("M13").Value = "Aleksandr Vlasov"
